# Insert two new weekly records (rows) right before the current row 313.
# The new rows represent a newer "Murcott" price report for the same
# market/quality combos already present at (old) rows 313-314; everything
# from (old) row 313 onward shifts down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 313 (old row 313 -> new row 315, etc.)
$ws.Rows.Item(313).Insert()
$ws.Rows.Item(313).Insert()

# Duplicate the (now shifted) original rows into the freshly inserted blanks
# so every column (style + value) starts out identical to its neighbour.
$ws.Range("A313:T313").Value = $ws.Range("A315:T315").Value2
$ws.Range("A314:T314").Value = $ws.Range("A316:T316").Value2

# Row 313: new Murcott / Primera report dated 44932 (was 44831 at old row 313)
$ws.Range("D313").Value = 44932
$ws.Range("M313").Value = 200
$ws.Range("N313").Value = 8000
$ws.Range("O313").Value = 8000
$ws.Range("P313").Value = 8000
$ws.Range("S313").Value = 444

# Row 314: new Murcott / Segunda report dated 44932 (was 44831 at old row 314)
$ws.Range("D314").Value = 44932
$ws.Range("M314").Value = 170
$ws.Range("N314").Value = 6000
$ws.Range("O314").Value = 6000
$ws.Range("P314").Value = 6000
$ws.Range("S314").Value = 333
